# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados..." timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 22:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1718902
$ws.Range("C4").Value = 12676
$ws.Range("D4").Value = 469558
$ws.Range("E4").Value = 1148968
$ws.Range("G4").Value = 571
$ws.Range("H4").Value = 100376

# Row 11 - Alemania
$ws.Range("B11").Value = 181288
$ws.Range("C11").Value = 499
$ws.Range("E11").Value = 10790
$ws.Range("G11").Value = 70
$ws.Range("H11").Value = 8498

# Row 13 - India
$ws.Range("B13").Value = 150793
$ws.Range("C13").Value = 5843
$ws.Range("E13").Value = 82167

# Row 34 - Sudafrica
$ws.Range("B34").Value = 24264
$ws.Range("C34").Value = 649
$ws.Range("D34").Value = 12741
$ws.Range("E34").Value = 10999
$ws.Range("G34").Value = 43
$ws.Range("H34").Value = 524

# Row 101 - Maldivas
$ws.Range("B101").Value = 1438
$ws.Range("C101").Value = 43
$ws.Range("D101").Value = 197
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 5
